$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing hyperlinks; will re-add them for the new row layout
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2025-11-09 01:53:31'
$ws.Range("B2").Value = '【AI開発】Microsoft Teams会議用AI BotのMVP開発依頼'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5429935'
$ws.Range("G2").Value = 458
$ws.Range("H2").Value = '🔥AI,Ai ★bot ◆開発'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5429935') | Out-Null
$ws.Range("F2").Style = "Hyperlink"

# Row 3
$ws.Range("A3").Value = '2025-11-09 01:53:31'
$ws.Range("B3").Value = '添付サイトのようなAIアバター生成iOSアプリ開発のエンジニアを募集します'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5429800'
$ws.Range("G3").Value = 388
$ws.Range("H3").Value = '🔥AI,Ai ◆開発 ◇アプリ'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5429800') | Out-Null
$ws.Range("F3").Style = "Hyperlink"

# Row 4
$ws.Range("A4").Value = '2025-11-09 01:53:31'
$ws.Range("B4").Value = 'WordPress保守業務および保守業務の自動化ツール開発'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5429668'
$ws.Range("G4").Value = 218
$ws.Range("H4").Value = '◆ツール,開発 ○WordPress'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5429668') | Out-Null
$ws.Range("F4").Style = "Hyperlink"

# Row 5
$ws.Range("A5").Value = '2025-11-09 01:53:31'
$ws.Range("B5").Value = '複数の見積書から情報抜き出して別表に書き出す作業の自動化 Excel VBAツール開発依頼'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5429304'
$ws.Range("G5").Value = 208
$ws.Range("H5").Value = '◆ツール,開発'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5429304') | Out-Null
$ws.Range("F5").Style = "Hyperlink"

# Row 6
$ws.Range("A6").Value = '2025-11-09 01:53:31'
$ws.Range("B6").Value = '【急募】ebayAPIを活用したShippingポリシー設定の専門家募集'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5415908'
$ws.Range("G6").Value = 183
$ws.Range("H6").Value = '🔥API'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5415908') | Out-Null
$ws.Range("F6").Style = "Hyperlink"

# Row 7
$ws.Range("A7").Value = '2025-11-09 01:53:31'
$ws.Range("B7").Value = '初回 自動売買ツール開発パートナー募集|IBKR(継続依頼あり)'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5429809'
$ws.Range("G7").Value = 135
$ws.Range("H7").Value = '◆ツール,開発'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5429809') | Out-Null
$ws.Range("F7").Style = "Hyperlink"

# Row 8
$ws.Range("A8").Value = '2025-11-09 01:53:31'
$ws.Range("B8").Value = 'Glideメインで作成したシステムをLinux+MySQL型に移行するための新規開発'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5430095'
$ws.Range("G8").Value = 115
$ws.Range("H8").Value = '◆開発 ◇MySQL'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5430095') | Out-Null
$ws.Range("F8").Style = "Hyperlink"

# Row 9
$ws.Range("A9").Value = '2025-11-09 01:53:31'
$ws.Range("B9").Value = '【急募】既存で作成済みのAccessデータベースをPower Apps等のアプリに移行したい'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5429495'
$ws.Range("G9").Value = 33
$ws.Range("H9").Value = '◇アプリ'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5429495') | Out-Null
$ws.Range("F9").Style = "Hyperlink"

# Row 10
$ws.Range("A10").Value = '2025-11-09 01:53:31'
$ws.Range("B10").Value = 'eBayテラピークでのキーワード検索結果等の取得するためのシステム制作'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5420779'
$ws.Range("G10").Value = 33
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5420779') | Out-Null
$ws.Range("F10").Style = "Hyperlink"

# Row 11
$ws.Range("A11").Value = '2025-11-09 01:53:31'
$ws.Range("B11").Value = '【フルスタックエンジニア】 働きながらスキルアップもできるEC業界で活躍してみませんか?'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5429335'
$ws.Range("G11").Value = 25
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5429335') | Out-Null
$ws.Range("F11").Style = "Hyperlink"

# Row 12
$ws.Range("A12").Value = '2025-11-09 01:53:31'
$ws.Range("B12").Value = 'UTAGE構築代行|ヒアリングから構築までお任せしたいです。'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5429882'
$ws.Range("G12").Value = 18
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5429882') | Out-Null
$ws.Range("F12").Style = "Hyperlink"

# Row 13
$ws.Range("A13").Value = '2025-11-09 01:53:31'
$ws.Range("B13").Value = 'MT4 RSXを使用したEAの作成依頼'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5430008'
$ws.Range("G13").Value = 10
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5430008') | Out-Null
$ws.Range("F13").Style = "Hyperlink"

# Row 14
$ws.Range("A14").Value = '2025-11-09 01:53:31'
$ws.Range("B14").Value = '【急募】LINE × QRコード連携で自動取得設定を実現!'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5430015'
$ws.Range("G14").Value = 10
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5430015') | Out-Null
$ws.Range("F14").Style = "Hyperlink"

$ws.Range("H10").ClearContents()
$ws.Range("H11").ClearContents()
$ws.Range("H12").ClearContents()
$ws.Range("H13").ClearContents()
$ws.Range("H14").ClearContents()

# Column D width: 28 -> 30 (input 29.15 compensates for the engine pixel-rounding on save)
$ws.Columns.Item(4).ColumnWidth = 29.15

$ws.Range("A1").Select() | Out-Null